$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.389.25'
$ws.Range('E2').Value = '  -1.78%  '
$ws.Range('D3').Value = '1.797.06'
$ws.Range('E3').Value = '  -1.53%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('E5').Value = '  -0.05%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '307.63'
$ws.Range('E6').Value = '  -0.97%  '
$ws.Range('E7').Value = '  -1.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3601'
$ws.Range('E8').Value = '  -2.26%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.98'
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07079'
$ws.Range('E10').Value = '  -0.97%  '
$ws.Range('E11').Value = '  +1.23%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07742'
$ws.Range('E12').Value = '  -0.05%  '
$ws.Range('E13').Value = '  -0.48%  '
$ws.Range('D14').Value = '1.799.08'
$ws.Range('E14').Value = '  -0.70%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '5.285'
$ws.Range('E15').Value = '  -0.47%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.329'
$ws.Range('E16').Value = '  -0.73%  '
$ws.Range('E17').Value = '  -1.88%  '
$ws.Range('E18').Value = '  -0.03%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.000008557'
$ws.Range('E19').Value = '  -1.74%  '
$ws.Range('B21').Value = 'Avalanche'
$ws.Range('C21').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.27'
$ws.Range('E21').Value = '  -1.20%  '
$ws.Range('B22').Value = 'WrappedBTC'
$ws.Range('C22').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D22').Value = '26.413.77'
$ws.Range('E22').Value = '  -1.82%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '4.974'
$ws.Range('E23').Value = '  -0.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.54'
$ws.Range('E24').Value = '  +1.11%  '
$ws.Range('D25').Value = '2.010.68'
$ws.Range('E25').Value = '  -2.03%  '
$ws.Range('E26').Value = '  -1.63%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '151.18'
$ws.Range('E28').Value = '  -1.85%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.025'
$ws.Range('E29').Value = '  +3.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '111.99'
$ws.Range('E30').Value = '  -1.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.858'
$ws.Range('E31').Value = '  -0.91%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.08691'
$ws.Range('E32').Value = '  -1.17%  '
$ws.Range('E33').Value = '  +1.03%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.747'
$ws.Range('E34').Value = '  +8.59%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '4.446'
$ws.Range('E35').Value = '  -0.85%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7239'
$ws.Range('E36').Value = '  -3.21%  '
$ws.Range('E37').Value = '  -2.40%  '
$ws.Range('E38').Value = '  +0.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.067'
$ws.Range('E39').Value = '  -1.52%  '
$ws.Range('E40').Value = '  -0.65%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.05100'
$ws.Range('E41').Value = '  -0.62%  '
$ws.Range('E42').Value = '  -1.89%  '
$ws.Range('B43').Value = 'TheSandbox'
$ws.Range('C43').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.5056'
$ws.Range('E43').Value = '  +1.87%  '
$ws.Range('B44').Value = 'FraxShare'
$ws.Range('C44').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '6.838'
$ws.Range('E44').Value = '  -1.13%  '
$ws.Range('E45').Value = '  -4.58%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.011'
$ws.Range('E46').Value = '  -3.57%  '
$ws.Range('E47').Value = '  +0.07%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.4627'
$ws.Range('E48').Value = '  -1.21%  '
$ws.Range('E49').Value = '  -2.60%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '100.92'
$ws.Range('E50').Value = '  -0.43%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.570'
$ws.Range('E51').Value = '  -2.46%  '
